$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 134, pushing the existing rows 134 and 135
# down to 135 and 136 respectively.
$ws.Rows("134:134").Insert()

# Populate the newly inserted row 134 with the new weekly data point.
$ws.Range("A134").Value = 4
$ws.Range("B134").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C134").Value = "Los Lagos"
$ws.Range("D134").Value = 44509
$ws.Range("D134").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E134").Value = 10
$ws.Range("F134").Value = 100112028
$ws.Range("G134").Value = "Sandia"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 600
$ws.Range("K134").Value = 1200
$ws.Range("L134").Value = 1200
$ws.Range("M134").Value = 1200
$ws.Range("N134").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O134").Value = "Perú"
$ws.Range("P134").Value = 1200
$ws.Range("Q134").Value = 1
$ws.Range("R134").Value = "Hortaliza"
